$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24, shifting existing rows 24-33 down to 25-34.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new record.
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44447
$ws.Range("D24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 100112022
$ws.Range("G24").Value = "Arveja Verde"
$ws.Range("H24").Value = "Perfection"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 35000
$ws.Range("L24").Value = 35000
$ws.Range("M24").Value = 35000
$ws.Range("N24").Value = "$/malla 25 kilos"
$ws.Range("O24").Value = "Provincia del Elquí"
$ws.Range("P24").Value = 1400
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
